$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual data cells per the diff
$ws.Range("G2").Value = 2
$ws.Range("G7").Value = 8
$ws.Range("G8").Value = 4
$ws.Range("G9").Value = 4
$ws.Range("G11").Value = 6
$ws.Range("G13").Value = 6
$ws.Range("D15").Value = 6

# New formatted (but empty) row 20, matching the style already used by the
# existing data rows (wrap text), which materializes the cells B20:J20.
$ws.Range("B20:J20").WrapText = $true

# Selection changes: active cell A7, with the whole row 7 selected
$ws.Range("A7:XFD7").Select()
